$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-number formatting from an existing date cell (A15) so the
# new date cells reuse the same cell style (numFmtId 14) instead of Excel
# creating a brand-new custom number format entry.
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 16 - UAE Tour, Stage 3
$ws.Range("A16").Value = 46071
$ws.Range("B16").Value = "UAE Tour"
$ws.Range("C16").Value = "Stage 3"
$ws.Range("D16").Value = "Antonio Tiberi"
$ws.Range("E16").Value = "Isaac del Toro"
$ws.Range("F16").Value = "Lennert Van Eetvelt"
$ws.Range("G16").Value = "Harold Tejada"
$ws.Range("H16").Value = "Felix Gall"
$ws.Range("I16").Value = "Tobias Halland Johannessen"
$ws.Range("J16").Value = "Derek Gee-West"
$ws.Range("K16").Value = "Jørgen Nordhagen"
$ws.Range("L16").Value = "Ilan Van Wilder"
$ws.Range("M16").Value = "Luke Plapp"

# Row 17 - Volta ao Algarve, Stage 1
$ws.Range("A17").Value = 46071
$ws.Range("B17").Value = "Volta ao Algarve"
$ws.Range("C17").Value = "Stage 1"
$ws.Range("D17").Value = "Paul Magnier"
$ws.Range("E17").Value = "Jordi Meeus"
$ws.Range("F17").Value = "Pavel Bittner"
$ws.Range("G17").Value = "Jasper Philipsen"
$ws.Range("H17").Value = "Santiago Mesa"
$ws.Range("I17").Value = "Arnaud De Lie"
$ws.Range("J17").Value = "Pascal Ackermann"
$ws.Range("K17").Value = "Kaden Groves"
$ws.Range("L17").Value = "John Degenkolb"
$ws.Range("M17").Value = "Oded Kogut"

# Selection / view tweaks to mirror the saved workbook state.
$null = $ws.Range("B21").Select()
